$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "Expected output" / "Result" cells that were left blank
$ws.Range("F24").Value = "It gets displayed as expected"
$ws.Range("G24").Value = "Pass"
$ws.Range("G38").Value = "Pass"
$ws.Range("G39").Value = "Pass"

# Scroll the sheet view so column D is the left-most visible column,
# then set the active selection to F25
$ws.Activate()
$excel.Goto($ws.Range("D1"), $true)
$ws.Range("F25").Select()
